$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.032687173732303
$ws.Cells.Item(2, 4).Value = 1.039950271089563
$ws.Cells.Item(2, 5).Value = 1.032093291896511
$ws.Cells.Item(2, 6).Value = 1.047209856671653
$ws.Cells.Item(2, 9).Value = 1.031709063085748
$ws.Cells.Item(2, 10).Value = 1.037816024603923
$ws.Cells.Item(2, 11).Value = 1.042734004010001
$ws.Cells.Item(2, 12).Value = 1.034899517529054
$ws.Cells.Item(2, 13).Value = 1.049973132548261
$ws.Cells.Item(2, 14).Value = 1.039289842398022
# Row 3
$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.033676102838914
$ws.Cells.Item(3, 4).Value = 1.040860434702759
$ws.Cells.Item(3, 5).Value = 1.0329346223286
$ws.Cells.Item(3, 6).Value = 1.048279545029488
$ws.Cells.Item(3, 9).Value = 1.031816832146939
$ws.Cells.Item(3, 10).Value = 1.03844715264282
$ws.Cells.Item(3, 11).Value = 1.043454228186014
$ws.Cells.Item(3, 12).Value = 1.035549460993314
$ws.Cells.Item(3, 13).Value = 1.050853950114456
$ws.Cells.Item(3, 14).Value = 1.039921866711126
# Row 4
$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.034316291881809
$ws.Cells.Item(4, 4).Value = 1.041449985217088
$ws.Cells.Item(4, 5).Value = 1.033479635549381
$ws.Cells.Item(4, 6).Value = 1.048972653971716
$ws.Cells.Item(4, 9).Value = 1.031885072854658
$ws.Cells.Item(4, 10).Value = 1.038855229477799
$ws.Cells.Item(4, 11).Value = 1.043920245143445
$ws.Cells.Item(4, 12).Value = 1.035969990744482
$ws.Cells.Item(4, 13).Value = 1.051424248843068
$ws.Cells.Item(4, 14).Value = 1.040330523062018
# Row 5
$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.034585495173064
$ws.Cells.Item(5, 4).Value = 1.041697978214467
$ws.Cells.Item(5, 5).Value = 1.033708905537454
$ws.Cells.Item(5, 6).Value = 1.049264263320927
$ws.Cells.Item(5, 9).Value = 1.031913403339089
$ws.Cells.Item(5, 10).Value = 1.039026711149956
$ws.Cells.Item(5, 11).Value = 1.044116153982059
$ws.Cells.Item(5, 12).Value = 1.036146774128747
$ws.Cells.Item(5, 13).Value = 1.051664085606877
$ws.Cells.Item(5, 14).Value = 1.040502248257818
# Row 6
$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.03463069954974
$ws.Cells.Item(6, 4).Value = 1.04173962586064
$ws.Cells.Item(6, 5).Value = 1.033747409532534
$ws.Cells.Item(6, 6).Value = 1.049313239062323
$ws.Cells.Item(6, 9).Value = 1.031918139150681
$ws.Cells.Item(6, 10).Value = 1.039055499325756
$ws.Cells.Item(6, 11).Value = 1.044149047639793
$ws.Cells.Item(6, 12).Value = 1.036176456383356
$ws.Cells.Item(6, 13).Value = 1.051704360108699
$ws.Cells.Item(6, 14).Value = 1.04053107731613
# Row 7
$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.034319888723629
$ws.Cells.Item(7, 4).Value = 1.041453298339361
$ws.Cells.Item(7, 5).Value = 1.033482698492174
$ws.Cells.Item(7, 6).Value = 1.048976549581777
$ws.Cells.Item(7, 9).Value = 1.031885452814893
$ws.Cells.Item(7, 10).Value = 1.03885752111535
$ws.Cells.Item(7, 11).Value = 1.043922862906369
$ws.Cells.Item(7, 12).Value = 1.03597235296253
$ws.Cells.Item(7, 13).Value = 1.051427453227836
$ws.Cells.Item(7, 14).Value = 1.040332817953957
# Row 8
$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.033021327419472
$ws.Cells.Item(8, 4).Value = 1.040257737462545
$ws.Cells.Item(8, 5).Value = 1.032377495221473
$ws.Cells.Item(8, 6).Value = 1.047571165840023
$ws.Cells.Item(8, 9).Value = 1.031745793018936
$ws.Cells.Item(8, 10).Value = 1.038029380017026
$ws.Cells.Item(8, 11).Value = 1.042977410023959
$ws.Cells.Item(8, 12).Value = 1.035119174085684
$ws.Cells.Item(8, 13).Value = 1.050270735826774
$ws.Cells.Item(8, 14).Value = 1.039503500800284
# Row 9
$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.030735311528638
$ws.Cells.Item(9, 4).Value = 1.038155750631575
$ws.Cells.Item(9, 5).Value = 1.030434752675461
$ws.Cells.Item(9, 6).Value = 1.045102006246605
$ws.Cells.Item(9, 9).Value = 1.031488279750389
$ws.Cells.Item(9, 10).Value = 1.036567785546345
$ws.Cells.Item(9, 11).Value = 1.041311312073088
$ws.Cells.Item(9, 12).Value = 1.033615592527348
$ws.Cells.Item(9, 13).Value = 1.048235175218837
$ws.Cells.Item(9, 14).Value = 1.038039830697808
# Row 10
$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.029212823611095
$ws.Cells.Item(10, 4).Value = 1.03675767270709
$ws.Cells.Item(10, 5).Value = 1.029142859046297
$ws.Cells.Item(10, 6).Value = 1.043460860945549
$ws.Cells.Item(10, 9).Value = 1.031308956550585
$ws.Cells.Item(10, 10).Value = 1.035591884084932
$ws.Cells.Item(10, 11).Value = 1.040200566388617
$ws.Cells.Item(10, 12).Value = 1.032613137619555
$ws.Cells.Item(10, 13).Value = 1.046880012832919
$ws.Cells.Item(10, 14).Value = 1.037062543344382
# Row 11
$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.028553937072561
$ws.Cells.Item(11, 4).Value = 1.036153071137967
$ws.Cells.Item(11, 5).Value = 1.028584242928323
$ws.Cells.Item(11, 6).Value = 1.042751413352071
$ws.Cells.Item(11, 9).Value = 1.031229499241063
$ws.Cells.Item(11, 10).Value = 1.035168959858299
$ws.Cells.Item(11, 11).Value = 1.039719607325865
$ws.Cells.Item(11, 12).Value = 1.03217905741032
$ws.Cells.Item(11, 13).Value = 1.046293667427991
$ws.Cells.Item(11, 14).Value = 1.036639018516837
# Row 12
$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.028309251909115
$ws.Cells.Item(12, 4).Value = 1.035928612569675
$ws.Cells.Item(12, 5).Value = 1.028376866571833
$ws.Cells.Item(12, 6).Value = 1.04248807080221
$ws.Cells.Item(12, 9).Value = 1.031199713861154
$ws.Cells.Item(12, 10).Value = 1.035011814680674
$ws.Cells.Item(12, 11).Value = 1.0395409585868
$ws.Cells.Item(12, 12).Value = 1.032017819840399
$ws.Cells.Item(12, 13).Value = 1.046075940771692
$ws.Cells.Item(12, 14).Value = 1.036481650175035
# Row 13
$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.028361735228013
$ws.Cells.Item(13, 4).Value = 1.035976754368948
$ws.Cells.Item(13, 5).Value = 1.028421344123875
$ws.Cells.Item(13, 6).Value = 1.042544550612301
$ws.Cells.Item(13, 9).Value = 1.031206115200944
$ws.Cells.Item(13, 10).Value = 1.035045525210257
$ws.Cells.Item(13, 11).Value = 1.039579279297093
$ws.Cells.Item(13, 12).Value = 1.032052405873147
$ws.Cells.Item(13, 13).Value = 1.046122640784095
$ws.Cells.Item(13, 14).Value = 1.036515408577436
# Row 14
$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.028533710203596
$ws.Cells.Item(14, 4).Value = 1.036134514923199
$ws.Cells.Item(14, 5).Value = 1.028567098709981
$ws.Cells.Item(14, 6).Value = 1.042729641741986
$ws.Cells.Item(14, 9).Value = 1.031227042706309
$ws.Cells.Item(14, 10).Value = 1.03515597125174
$ws.Cells.Item(14, 11).Value = 1.039704840143537
$ws.Cells.Item(14, 12).Value = 1.032165729475593
$ws.Cells.Item(14, 13).Value = 1.046275668680693
$ws.Cells.Item(14, 14).Value = 1.036626011464967
# Row 15
$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.028639676914439
$ws.Cells.Item(15, 4).Value = 1.036231731986999
$ws.Cells.Item(15, 5).Value = 1.028656918658967
$ws.Cells.Item(15, 6).Value = 1.042843706087462
$ws.Cells.Item(15, 9).Value = 1.031239900882599
$ws.Cells.Item(15, 10).Value = 1.035224013794105
$ws.Cells.Item(15, 11).Value = 1.03978220246148
$ws.Cells.Item(15, 12).Value = 1.032235551796426
$ws.Cells.Item(15, 13).Value = 1.046369963267736
$ws.Cells.Item(15, 14).Value = 1.036694150635547
# Row 16
$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.029256559378535
$ws.Cells.Item(16, 4).Value = 1.03679781453383
$ws.Cells.Item(16, 5).Value = 1.029179949127107
$ws.Cells.Item(16, 6).Value = 1.043507969550606
$ws.Cells.Item(16, 9).Value = 1.031314191788306
$ws.Cells.Item(16, 10).Value = 1.035619944803768
$ws.Cells.Item(16, 11).Value = 1.040232486141448
$ws.Cells.Item(16, 12).Value = 1.032641945933086
$ws.Cells.Item(16, 13).Value = 1.04691893616278
$ws.Cells.Item(16, 14).Value = 1.037090643912657
# Row 17
$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.029643610343927
$ws.Cells.Item(17, 4).Value = 1.037153111376585
$ws.Cells.Item(17, 5).Value = 1.029508242486766
$ws.Cells.Item(17, 6).Value = 1.043924960671522
$ws.Cells.Item(17, 9).Value = 1.031360308432246
$ws.Cells.Item(17, 10).Value = 1.035868207971798
$ws.Cells.Item(17, 11).Value = 1.040514937998235
$ws.Cells.Item(17, 12).Value = 1.032896863919474
$ws.Cells.Item(17, 13).Value = 1.047263413019432
$ws.Cells.Item(17, 14).Value = 1.037339259642862
# Row 18
$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.029869405404235
$ws.Cells.Item(18, 4).Value = 1.037360424769456
$ws.Cells.Item(18, 5).Value = 1.029699805867231
$ws.Cells.Item(18, 6).Value = 1.044168298260008
$ws.Cells.Item(18, 9).Value = 1.031387032903113
$ws.Cells.Item(18, 10).Value = 1.03601298162623
$ws.Cells.Item(18, 11).Value = 1.040679687451461
$ws.Cells.Item(18, 12).Value = 1.033045552193676
$ws.Cells.Item(18, 13).Value = 1.047464383802413
$ws.Cells.Item(18, 14).Value = 1.037484238892489
# Row 19
$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.029946401571061
$ws.Cells.Item(19, 4).Value = 1.037431125932076
$ws.Cells.Item(19, 5).Value = 1.029765136791217
$ws.Cells.Item(19, 6).Value = 1.04425128935673
$ws.Cells.Item(19, 9).Value = 1.0313961156315
$ws.Cells.Item(19, 10).Value = 1.036062339905272
$ws.Cells.Item(19, 11).Value = 1.040735862747088
$ws.Cells.Item(19, 12).Value = 1.033096250847908
$ws.Cells.Item(19, 13).Value = 1.047532916962362
$ws.Cells.Item(19, 14).Value = 1.037533667265949
# Row 20
$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.029602079837641
$ws.Cells.Item(20, 4).Value = 1.037114983649119
$ws.Cells.Item(20, 5).Value = 1.029473011908078
$ws.Cells.Item(20, 6).Value = 1.04388020970248
$ws.Cells.Item(20, 9).Value = 1.031355378608191
$ws.Cells.Item(20, 10).Value = 1.035841575200645
$ws.Cells.Item(20, 11).Value = 1.040484633575393
$ws.Cells.Item(20, 12).Value = 1.032869513737775
$ws.Cells.Item(20, 13).Value = 1.047226449433053
$ws.Cells.Item(20, 14).Value = 1.037312589050119
# Row 21
$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.028483066348004
$ws.Cells.Item(21, 4).Value = 1.036088055124875
$ws.Cells.Item(21, 5).Value = 1.028524174336815
$ws.Cells.Item(21, 6).Value = 1.042675132101562
$ws.Cells.Item(21, 9).Value = 1.031220887566171
$ws.Cells.Item(21, 10).Value = 1.035123449083426
$ws.Cells.Item(21, 11).Value = 1.039667865574103
$ws.Cells.Item(21, 12).Value = 1.032132358516766
$ws.Cells.Item(21, 13).Value = 1.046230603892802
$ws.Cells.Item(21, 14).Value = 1.036593443111443
# Row 22
$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.027779814464966
$ws.Cells.Item(22, 4).Value = 1.035443063931035
$ws.Cells.Item(22, 5).Value = 1.027928288840879
$ws.Cells.Item(22, 6).Value = 1.041918481470362
$ws.Cells.Item(22, 9).Value = 1.031134757692854
$ws.Cells.Item(22, 10).Value = 1.034671632491132
$ws.Cells.Item(22, 11).Value = 1.039154336503643
$ws.Cells.Item(22, 12).Value = 1.03166887520553
$ws.Cells.Item(22, 13).Value = 1.045604870827522
$ws.Cells.Item(22, 14).Value = 1.036140984887765
# Row 23
$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.0281525914664
$ws.Cells.Item(23, 4).Value = 1.035784921261196
$ws.Cells.Item(23, 5).Value = 1.028244113624471
$ws.Cells.Item(23, 6).Value = 1.04231949846162
$ws.Cells.Item(23, 9).Value = 1.031180565450342
$ws.Cells.Item(23, 10).Value = 1.034911177438845
$ws.Cells.Item(23, 11).Value = 1.039426567235654
$ws.Cells.Item(23, 12).Value = 1.031914576663034
$ws.Cells.Item(23, 13).Value = 1.045936546131749
$ws.Cells.Item(23, 14).Value = 1.036380870016779
# Row 24
$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.029620845578252
$ws.Cells.Item(24, 4).Value = 1.037132211697025
$ws.Cells.Item(24, 5).Value = 1.029488930857141
$ws.Cells.Item(24, 6).Value = 1.04390043038772
$ws.Cells.Item(24, 9).Value = 1.031357606722951
$ws.Cells.Item(24, 10).Value = 1.035853609508762
$ws.Cells.Item(24, 11).Value = 1.040498326839456
$ws.Cells.Item(24, 12).Value = 1.032881872110962
$ws.Cells.Item(24, 13).Value = 1.047243151552994
$ws.Cells.Item(24, 14).Value = 1.037324640448334
# Row 25
$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.031326035508108
$ws.Cells.Item(25, 4).Value = 1.038698596539418
$ws.Cells.Item(25, 5).Value = 1.030936426790618
$ws.Cells.Item(25, 6).Value = 1.04573947234743
$ws.Cells.Item(25, 9).Value = 1.031556203038422
$ws.Cells.Item(25, 10).Value = 1.036945911263793
$ws.Cells.Item(25, 11).Value = 1.041742044135108
$ws.Cells.Item(25, 12).Value = 1.034004319760657
$ws.Cells.Item(25, 13).Value = 1.048761088866425
$ws.Cells.Item(25, 14).Value = 1.038418493397146
